$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 762
$ws1.Range("F6").Value = 2456
$ws1.Range("F8").Value = 1821
$ws1.Range("F9").Value = 3092
$ws1.Range("F10").Value = 188
$ws1.Range("F11").Value = 4588
$ws1.Range("F12").Value = 418
$ws1.Range("F13").Value = 240
$ws1.Range("F15").Value = 586
$ws1.Range("F16").Value = 274
$ws1.Range("F17").Value = 626
$ws1.Range("F18").Value = 258
$ws1.Range("F19").Value = 8
$ws1.Range("F21").Value = 123
$ws1.Range("F22").Value = 318
$ws1.Range("F23").Value = 4593
$ws1.Range("F24").Value = 6
$ws1.Range("F25").Value = 18
$ws1.Range("F27").Value = 4412
$ws1.Range("F28").Value = 7
$ws1.Range("F29").Value = 1157
$ws1.Range("F31").Value = 610
$ws1.Range("F32").Value = 4390
$ws1.Range("F33").Value = 44
$ws1.Range("F35").Value = 703
$ws1.Range("F36").Value = 30
$ws1.Range("F37").Value = 646
$ws1.Range("F38").Value = 638

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 5

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1056
$ws3.Range("F4").Value = 26

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1056
$ws4.Range("F5").Value = 26
$ws4.Range("F8").Value = 762
$ws4.Range("F9").Value = 2456
$ws4.Range("F11").Value = 1821
$ws4.Range("F13").Value = 3092
$ws4.Range("F14").Value = 188
$ws4.Range("F15").Value = 4588
$ws4.Range("F16").Value = 418
$ws4.Range("F17").Value = 240
$ws4.Range("F19").Value = 586
$ws4.Range("F20").Value = 274
$ws4.Range("F21").Value = 626
$ws4.Range("F22").Value = 258
$ws4.Range("F23").Value = 8
$ws4.Range("F26").Value = 123
$ws4.Range("F27").Value = 318
$ws4.Range("F28").Value = 4593
$ws4.Range("F29").Value = 6
$ws4.Range("F30").Value = 18
$ws4.Range("F32").Value = 4412
$ws4.Range("F33").Value = 7
$ws4.Range("F34").Value = 1157
$ws4.Range("F36").Value = 610
$ws4.Range("F37").Value = 4390
$ws4.Range("F38").Value = 5
$ws4.Range("F39").Value = 44
$ws4.Range("F41").Value = 703
$ws4.Range("F42").Value = 30
$ws4.Range("F43").Value = 646
$ws4.Range("F44").Value = 638
